# 4.0.3 model and data
# Expand the single "BVTQaZ" and "VTQaZ" trans CSV rows on the "Boolean"
# sheet into their per-vehicle-type breakdowns (LDVs/HDVs/aircraft/rail/
# ships/motorbikes), matching the new InputData file layout.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("About")
$ws2 = $wb.Worksheets.Item("Integer")
$ws3 = $wb.Worksheets.Item("Boolean")

# --- "Boolean" sheet: split trans/BVTQaZ/BVTQaZ.csv and trans/VTQaZ/VTQaZ.csv ---
# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv". Insert 5 more rows right
# below it (keeping row 17 itself untouched for now) so the six sub-files can
# be written in order, then overwrite rows 17-22.
$ws3.Range("A18:A22").Insert()

$ws3.Cells.Item(17, 1).Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$ws3.Cells.Item(18, 1).Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$ws3.Cells.Item(19, 1).Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$ws3.Cells.Item(20, 1).Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$ws3.Cells.Item(21, 1).Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$ws3.Cells.Item(22, 1).Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# After the insert above, "trans/BVTStL/BVTStL.csv" / PVTStL / SRPbVT sit at
# rows 23-25, and "trans/VTQaZ/VTQaZ.csv" now sits at row 26, followed by
# "trans/VTStFES/VTStFES.csv" at row 27. Repeat the same expansion there.
$ws3.Range("A27:A31").Insert()

$ws3.Cells.Item(26, 1).Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$ws3.Cells.Item(27, 1).Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$ws3.Cells.Item(28, 1).Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$ws3.Cells.Item(29, 1).Value = "trans/VTQaZ/VTQaZ-rail.csv"
$ws3.Cells.Item(30, 1).Value = "trans/VTQaZ/VTQaZ-ships.csv"
$ws3.Cells.Item(31, 1).Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# A handful of trailing rows (33-38) carry formatting only (no values), as
# in the source workbook.
for ($r = 33; $r -le 38; $r++) {
    $cell = $ws3.Cells.Item($r, 1)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.ClearContents()
}

# --- View state: restore per-sheet selections and make "About" the active tab ---
$ws2.Activate()
$ws2.Range("A13").Select() | Out-Null

$ws3.Activate()
$ws3.Range("A32").Select() | Out-Null

$ws1.Activate()
